$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a title/label in the previously-empty row 7, just above the table (no row shift)
$ws.Range("D7").Value = "Déplacement des pions"

# Update the selected cell in the sheet view to E5
$ws.Range("E5").Select()
